$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes
$ws.Columns.Item(4).ColumnWidth = 49.166666666666664
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 44.166666666666664
$ws.Columns.Item(9).ColumnWidth = 38.166666666666664
$ws.Columns.Item(10).ColumnWidth = 37.166666666666664
$ws.Columns.Item(11).ColumnWidth = 36.166666666666664
$ws.Columns.Item(12).ColumnWidth = 44.166666666666664
$ws.Columns.Item(13).ColumnWidth = 36.166666666666664

# Cell value changes
$ws.Range('F2').Value = 'studentorganiza@cookman.edu'
$ws.Range('H2').Value = 'https://linkedin.com/groups/studentorganiza'
$ws.Range('K2').Value = 'https://twitter.com/studentorganiza'
$ws.Range('F3').Value = 'greekletteredor@cookman.edu'
$ws.Range('G3').Value = '(555) 456-7890'
$ws.Range('H3').Value = 'https://linkedin.com/groups/greekletteredor'
$ws.Range('J3').Value = 'https://facebook.com/greekletteredor'
$ws.Range('F4').Value = 'greeklifeinfofo@cookman.edu'
$ws.Range('H4').Value = 'https://linkedin.com/groups/greeklifeinfofo'
$ws.Range('I4').Value = 'https://instagram.com/greeklifeinfofo'
$ws.Range('D5').Value = 'https://cookman.edu/logos/studentexperien_logo.png'
$ws.Range('E5').Value = 'Student organization focused on general activities and community engagement. The STUDENT EXPERIENCE welcomes all interested students to participate and make a positive impact.'
$ws.Range('G5').Value = '(555) 456-7890'
$ws.Range('J5').Value = 'https://facebook.com/studentexperien'
$ws.Range('K5').Value = 'https://twitter.com/studentexperien'
$ws.Range('E6').Value = 'Student organization focused on general activities and community engagement. The WHO WE ARE welcomes all interested students to participate and make a positive impact.'
$ws.Range('F6').Value = 'whoweare@cookman.edu'
$ws.Range('D7').Value = 'https://cookman.edu/logos/wildcatweb_logo.png'
$ws.Range('E7').Value = 'Student organization focused on general activities and community engagement. The WILDCAT WEB welcomes all interested students to participate and make a positive impact.'
$ws.Range('G7').Value = '(555) 901-2345'
$ws.Range('I7').Value = 'https://instagram.com/wildcatweb'
$ws.Range('K7').Value = 'https://twitter.com/wildcatweb'
$ws.Range('E8').Value = 'Student organization focused on general activities and community engagement. The Financial Reports welcomes all interested students to participate and make a positive impact.'
$ws.Range('G8').Value = '(555) 234-5678'
$ws.Range('H8').Value = 'https://linkedin.com/groups/financialreport'
$ws.Range('J8').Value = 'https://facebook.com/financialreport'
$ws.Range('K8').Value = 'https://twitter.com/financialreport'
$ws.Range('M8').Value = 'https://tiktok.com/@financialreport'
$ws.Range('D9').Value = 'https://cookman.edu/logos/sacscocriseqep_logo.png'
$ws.Range('E9').Value = 'Student organization focused on general activities and community engagement. The SACSCOC RISE QEP welcomes all interested students to participate and make a positive impact.'
$ws.Range('G9').Value = '(555) 567-8901'
$ws.Range('K9').Value = 'https://twitter.com/sacscocriseqep'
$ws.Range('E10').Value = 'Student organization focused on general activities and community engagement. The Marketing and Communications welcomes all interested students to participate and make a positive impact.'
$ws.Range('F10').Value = 'marketingandcom@cookman.edu'
$ws.Range('G10').Value = '(555) 123-4567'
$ws.Range('K10').Value = 'https://twitter.com/marketingandcom'
$ws.Range('E11').Value = 'Student organization focused on general activities and community engagement. The B-CU Jobs welcomes all interested students to participate and make a positive impact.'
$ws.Range('F11').Value = 'bcujobs@cookman.edu'
$ws.Range('G11').Value = '(555) 123-4567'
$ws.Range('J11').Value = 'https://facebook.com/bcujobs'
$ws.Range('K11').Value = 'https://twitter.com/bcujobs'
$ws.Range('L11').Value = 'https://youtube.com/channel/bcujobs'
$ws.Range('E12').Value = 'Student organization focused on general activities and community engagement. The Payment Center welcomes all interested students to participate and make a positive impact.'
$ws.Range('F12').Value = 'paymentcenter@cookman.edu'
$ws.Range('M12').Value = 'https://tiktok.com/@paymentcenter'
$ws.Range('E13').Value = 'Student organization focused on general activities and community engagement. The Clery Report welcomes all interested students to participate and make a positive impact.'
$ws.Range('F13').Value = 'cleryreport@cookman.edu'
$ws.Range('G13').Value = '(555) 456-7890'
$ws.Range('I13').Value = 'https://instagram.com/cleryreport'
$ws.Range('L13').Value = 'https://youtube.com/channel/cleryreport'
$ws.Range('E14').Value = 'Student organization focused on general activities and community engagement. The Accessibility Statement welcomes all interested students to participate and make a positive impact.'
$ws.Range('F14').Value = 'accessibilityst@cookman.edu'
$ws.Range('G14').Value = '(555) 567-8901'
$ws.Range('K14').Value = 'https://twitter.com/accessibilityst'
$ws.Range('E15').Value = 'Student organization focused on general activities and community engagement. The Open Bids welcomes all interested students to participate and make a positive impact.'
$ws.Range('F15').Value = 'openbids@cookman.edu'
$ws.Range('H15').Value = 'https://linkedin.com/groups/openbids'
$ws.Range('J15').Value = 'https://facebook.com/openbids'
$ws.Range('E16').Value = 'Student organization focused on general activities and community engagement. The Residence Life welcomes all interested students to participate and make a positive impact.'
$ws.Range('F16').Value = 'residencelife@cookman.edu'
$ws.Range('G16').Value = '(555) 678-9012'
$ws.Range('L16').Value = 'https://youtube.com/channel/residencelife'
$ws.Range('D17').Value = 'https://cookman.edu/logos/centerforcivice_logo.png'
$ws.Range('E17').Value = 'Student organization focused on general activities and community engagement. The Center for Civic Engagement welcomes all interested students to participate and make a positive impact.'
$ws.Range('H17').Value = 'https://linkedin.com/groups/centerforcivice'
$ws.Range('J17').Value = 'https://facebook.com/centerforcivice'
$ws.Range('K17').Value = 'https://twitter.com/centerforcivice'
$ws.Range('E18').Value = 'Interfaith organization promoting spiritual growth and religious dialogue. The Chaplaincy & Religious Life welcomes all interested students to participate and make a positive impact.'
$ws.Range('I18').Value = 'https://instagram.com/chaplaincyrelig'
$ws.Range('K18').Value = 'https://twitter.com/chaplaincyrelig'
$ws.Range('L18').Value = 'https://youtube.com/channel/chaplaincyrelig'
$ws.Range('E19').Value = 'Student organization focused on general activities and community engagement. The Future Students welcomes all interested students to participate and make a positive impact.'
$ws.Range('F19').Value = 'futurestudents@cookman.edu'
$ws.Range('I19').Value = 'https://instagram.com/futurestudents'
$ws.Range('J19').Value = 'https://facebook.com/futurestudents'
$ws.Range('K19').Value = 'https://twitter.com/futurestudents'
$ws.Range('E20').Value = 'Student organization focused on general activities and community engagement. The Current Students welcomes all interested students to participate and make a positive impact.'
$ws.Range('F20').Value = 'currentstudents@cookman.edu'
$ws.Range('G20').Value = '(555) 890-1234'
$ws.Range('H20').Value = 'https://linkedin.com/groups/currentstudents'
$ws.Range('J20').Value = 'https://facebook.com/currentstudents'
$ws.Range('D21').Value = 'https://cookman.edu/logos/studentorganiza_logo.png'
$ws.Range('F21').Value = 'studentorganiza@cookman.edu'
$ws.Range('G21').Value = '(555) 234-5678'
$ws.Range('K21').Value = 'https://twitter.com/studentorganiza'
